# The underlying OOXML diff swaps the content of ppt/theme/theme1.xml
# (the presentation/slide-master theme, originally "Integral" / "Red
# Violet") with the content that used to live in ppt/theme/theme2.xml
# (the notes-master theme, "Office Theme" / "Office"), and vice versa.
#
# The PowerPoint object model's writable surface for theme content is
# the 12-slot DrawingML colour scheme reachable from a slide's
# ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -
# it always maps onto the single presentation theme part
# (ppt/theme/theme1.xml). Re-point each of those 12 slots at the
# "Office Theme" palette that the target state puts into theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (slot, target RGB) using the same ordering PowerPoint uses:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}

Write-Output "theme colours updated"
